# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the 65ecf2ef-1163-42b9-8d98-1526882d18ba.md row (row 4) that was
# just handed back, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 65ecf2ef-... row.
$overview.Range("G4").Value = "2016-10-24 07:18:51"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback
# DateTime" (K) for the same row.
$zhcn.Range("H4").Value = "2016-10-24 07:18:40"
$zhcn.Range("K4").Value = "2016-10-24 07:19:20"

# de-de sheet: same two columns.
$dede.Range("H4").Value = "2016-10-24 07:18:51"
$dede.Range("K4").Value = "2016-10-24 07:19:37"
